$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("D-P")

# Update the raw input values that drive the J7:L9 / J26:L28 calculation blocks.
$ws.Range("E17").Value = 6.2991999999999999
$ws.Range("F17").Value = -2
$ws.Range("F18").Value = 4
$ws.Range("G19").Value = 5

# Recalculate the workbook so all dependent formulas (including chart caches)
# pick up the new values.
$excel.CalculateFullRebuild()
$wb.RefreshAll()

# Restore the active sheet/cell selection recorded in the saved file.
$ws.Activate()
$ws.Range("AB4").Select()
